$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Activities sheet: insert a new "Activity Status" column between
#    Activity Name (B) and BL Project Finish (C), shifting the remaining
#    columns right.
# ---------------------------------------------------------------------------
$activities = $wb.Worksheets.Item("Activities")
$activities.Columns("C").Insert()

$activities.Range("C1").Value = "Activity Status"
$activities.Range("C4").Value = "Completed"
$activities.Range("C5").Value = "In Progress"
$activities.Range("C6").Value = "Not Started"

# ---------------------------------------------------------------------------
# 2. Resource assignment sheets: indent the Activity ID values to reflect
#    the same hierarchy used on the Activities sheet.
# ---------------------------------------------------------------------------
$resourceSheets = @(
    "Ressource Assign. Budgeted",
    "Ressource Assign. Actual",
    "Ressource Assign. Remaining"
)

foreach ($sheetName in $resourceSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A3").Value = "  A-110"
    $ws.Range("A4").Value = "    A-111"
    $ws.Range("A5").Value = "    A-112"
    $ws.Range("A6").Value = "  A-120"
}

# ---------------------------------------------------------------------------
# 3. README sheet: widen column A, add the two new guidance lines and
#    rewrite the leaf-activity note.
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("README")
$readme.Columns("A").ColumnWidth = 89.16666666666667

# Insert a row after row 6 for the new "Activity Status" note.
$readme.Rows("7").Insert()
$readme.Range("A6").Value = "- Activity Name is required only for leaf activities."
$readme.Range("A7").Value = "- Activity Status must be filled when Activity Name is filled."

# Insert a row after row 11 for the new indentation note.
$readme.Rows("12").Insert()
$readme.Range("A12").Value = "- Keep the same indentation in Activity ID for all tables."
